$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

# ---- Part 1: in-place cell updates on existing rows (status + stat corrections) ----
# Row 13
$ws.Range("G13").Value = '3:36 - 2nd Half'
$ws.Range("K13").Value = 3
$ws.Range("P13").Value = 30
$ws.Range("R13").Value = 13
$ws.Range("T13").Value = 7

# Row 23
$ws.Range("G23").Value = '3:36 - 2nd Half'
$ws.Range("H23").Value = 21
$ws.Range("I23").Value = 22
$ws.Range("J23").Value = 5
$ws.Range("P23").Value = 35
$ws.Range("U23").Value = 2
$ws.Range("V23").Value = 3

# Row 24
$ws.Range("G24").Value = '3:36 - 2nd Half'

# Row 32
$ws.Range("G32").Value = '3:36 - 2nd Half'

# Row 33
$ws.Range("D33").Value = 'Taylor Bol Bowen'
$ws.Range("G33").Value = '3:36 - 2nd Half'
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = 2

# Row 34
$ws.Range("D34").Value = 'Jalil Bethea'
$ws.Range("G34").Value = '3:36 - 2nd Half'
$ws.Range("H34").Value = -1
$ws.Range("J34").Value = 1
$ws.Range("O34").Value = 1
$ws.Range("P34").Value = 4
$ws.Range("R34").Value = 2
$ws.Range("T34").Value = 2

# Row 44
$ws.Range("G44").Value = '3:36 - 2nd Half'
$ws.Range("H44").Value = 34
$ws.Range("I44").Value = 25
$ws.Range("J44").Value = 7
$ws.Range("K44").Value = 6
$ws.Range("P44").Value = 33
$ws.Range("Q44").Value = 11
$ws.Range("R44").Value = 15
$ws.Range("V44").Value = 7

# Row 45
$ws.Range("G45").Value = '3:36 - 2nd Half'
$ws.Range("H45").Value = 25
$ws.Range("I45").Value = 15
$ws.Range("J45").Value = 2
$ws.Range("K45").Value = 8
$ws.Range("P45").Value = 34
$ws.Range("U45").Value = 3
$ws.Range("V45").Value = 4

# Row 46
$ws.Range("G46").Value = '3:36 - 2nd Half'

# Row 56
$ws.Range("D56").Value = 'Xaivian Lee'
$ws.Range("E56").Value = 'FLA'
$ws.Range("G56").Value = '3:36 - 2nd Half'
$ws.Range("I56").Value = 6
$ws.Range("K56").Value = 2
$ws.Range("N56").Value = 0
$ws.Range("O56").Value = 3
$ws.Range("P56").Value = 17
$ws.Range("Q56").Value = 3
$ws.Range("R56").Value = 11
$ws.Range("S56").Value = 0
$ws.Range("T56").Value = 3

# Row 57
$ws.Range("D57").Value = 'Amari Allen'
$ws.Range("E57").Value = 'ALA'
$ws.Range("G57").Value = '3:36 - 2nd Half'
$ws.Range("H57").Value = 3
$ws.Range("I57").Value = 5
$ws.Range("K57").Value = 3
$ws.Range("N57").Value = 3
$ws.Range("O57").Value = 5
$ws.Range("P57").Value = 21
$ws.Range("Q57").Value = 2
$ws.Range("R57").Value = 8
$ws.Range("S57").Value = 1
$ws.Range("T57").Value = 5

# Row 66
$ws.Range("G66").Value = '3:36 - 2nd Half'
$ws.Range("H66").Value = 28
$ws.Range("I66").Value = 14
$ws.Range("J66").Value = 17
$ws.Range("P66").Value = 29
$ws.Range("Q66").Value = 4
$ws.Range("R66").Value = 12
$ws.Range("U66").Value = 6
$ws.Range("V66").Value = 7

# Row 67
$ws.Range("G67").Value = '3:36 - 2nd Half'
$ws.Range("H67").Value = 9
$ws.Range("I67").Value = 11
$ws.Range("N67").Value = 2
$ws.Range("O67").Value = 4
$ws.Range("P67").Value = 28
$ws.Range("Q67").Value = 3
$ws.Range("R67").Value = 7
$ws.Range("S67").Value = 3
$ws.Range("T67").Value = 7

# Row 79
$ws.Range("G79").Value = '3:36 - 2nd Half'
$ws.Range("H79").Value = 20
$ws.Range("K79").Value = 5
$ws.Range("P79").Value = 27

# Row 147
$ws.Range("G147").Value = '3:36 - 2nd Half'
$ws.Range("H147").Value = 14
$ws.Range("I147").Value = 11
$ws.Range("J147").Value = 4
$ws.Range("K147").Value = 4
$ws.Range("P147").Value = 27
$ws.Range("U147").Value = 2
$ws.Range("V147").Value = 3

# Row 148
$ws.Range("G148").Value = '3:36 - 2nd Half'
$ws.Range("H148").Value = 12
$ws.Range("J148").Value = 7
$ws.Range("O148").Value = 5
$ws.Range("P148").Value = 24

# Row 149
$ws.Range("D149").Value = 'London Jemison'
$ws.Range("E149").Value = 'ALA'
$ws.Range("G149").Value = '3:36 - 2nd Half'
$ws.Range("H149").Value = 5
$ws.Range("I149").Value = 6
$ws.Range("M149").Value = 1
$ws.Range("O149").Value = 1
$ws.Range("P149").Value = 12
$ws.Range("Q149").Value = 2
$ws.Range("R149").Value = 4
$ws.Range("S149").Value = 2
$ws.Range("T149").Value = 3
$ws.Range("U149").Value = 0
$ws.Range("V149").Value = 0

# Row 150
$ws.Range("D150").Value = 'Isaiah Brown'
$ws.Range("E150").Value = 'FLA'
$ws.Range("G150").Value = '3:36 - 2nd Half'
$ws.Range("M150").Value = 0
$ws.Range("O150").Value = 2
$ws.Range("P150").Value = 11
$ws.Range("R150").Value = 1
$ws.Range("S150").Value = 0
$ws.Range("T150").Value = 0
$ws.Range("U150").Value = 1
$ws.Range("V150").Value = 2

# ---- Part 2: insert a new row at 151 (pushes old row151 -> row152), then populate it ----
$ws.Rows.Item(151).Insert()

# Row 151 (new): Olivier Rioux
$ws.Range("A151").NumberFormat = "@"
$ws.Range("A151").Value = "2026-02-01"
$ws.Range("B151").Value = 'Undrafted'
$ws.Range("C151").Value = 'No'
$ws.Range("D151").Value = 'Olivier Rioux'
$ws.Range("E151").Value = 'FLA'
$ws.Range("F151").Value = 'ALA@FLA'
$ws.Range("G151").Value = '3:36 - 2nd Half'
$ws.Range("H151").Value = 2
$ws.Range("I151").Value = 2
$ws.Range("J151").Value = 1
$ws.Range("K151").Value = 0
$ws.Range("L151").Value = 0
$ws.Range("M151").Value = 0
$ws.Range("N151").Value = 0
$ws.Range("O151").Value = 0
$ws.Range("P151").Value = 1
$ws.Range("Q151").Value = 1
$ws.Range("R151").Value = 2
$ws.Range("S151").Value = 0
$ws.Range("T151").Value = 0
$ws.Range("U151").Value = 0
$ws.Range("V151").Value = 0

# Row 152 (shifted down from old row151, Noah Williamson): only status changes
$ws.Range("G152").Value = '3:36 - 2nd Half'

# ---- Part 3: append new rows 153-156 ----
# Row 153: Alex Kovatchev
$ws.Range("A153").NumberFormat = "@"
$ws.Range("A153").Value = "2026-02-01"
$ws.Range("B153").Value = 'Undrafted'
$ws.Range("C153").Value = 'No'
$ws.Range("D153").Value = 'Alex Kovatchev'
$ws.Range("E153").Value = 'FLA'
$ws.Range("F153").Value = 'ALA@FLA'
$ws.Range("G153").Value = '3:36 - 2nd Half'
$ws.Range("H153").Value = 0
$ws.Range("I153").Value = 0
$ws.Range("J153").Value = 0
$ws.Range("K153").Value = 0
$ws.Range("L153").Value = 0
$ws.Range("M153").Value = 0
$ws.Range("N153").Value = 0
$ws.Range("O153").Value = 0
$ws.Range("P153").Value = 1
$ws.Range("Q153").Value = 0
$ws.Range("R153").Value = 0
$ws.Range("S153").Value = 0
$ws.Range("T153").Value = 0
$ws.Range("U153").Value = 0
$ws.Range("V153").Value = 0

# Row 154: Alex Lloyd
$ws.Range("A154").NumberFormat = "@"
$ws.Range("A154").Value = "2026-02-01"
$ws.Range("B154").Value = 'Undrafted'
$ws.Range("C154").Value = 'No'
$ws.Range("D154").Value = 'Alex Lloyd'
$ws.Range("E154").Value = 'FLA'
$ws.Range("F154").Value = 'ALA@FLA'
$ws.Range("G154").Value = '3:36 - 2nd Half'
$ws.Range("H154").Value = 0
$ws.Range("I154").Value = 0
$ws.Range("J154").Value = 0
$ws.Range("K154").Value = 0
$ws.Range("L154").Value = 0
$ws.Range("M154").Value = 0
$ws.Range("N154").Value = 0
$ws.Range("O154").Value = 0
$ws.Range("P154").Value = 1
$ws.Range("Q154").Value = 0
$ws.Range("R154").Value = 0
$ws.Range("S154").Value = 0
$ws.Range("T154").Value = 0
$ws.Range("U154").Value = 0
$ws.Range("V154").Value = 0

# Row 155: CJ Ingram
$ws.Range("A155").NumberFormat = "@"
$ws.Range("A155").Value = "2026-02-01"
$ws.Range("B155").Value = 'Undrafted'
$ws.Range("C155").Value = 'No'
$ws.Range("D155").Value = 'CJ Ingram'
$ws.Range("E155").Value = 'FLA'
$ws.Range("F155").Value = 'ALA@FLA'
$ws.Range("G155").Value = '3:36 - 2nd Half'
$ws.Range("H155").Value = 0
$ws.Range("I155").Value = 0
$ws.Range("J155").Value = 0
$ws.Range("K155").Value = 0
$ws.Range("L155").Value = 0
$ws.Range("M155").Value = 0
$ws.Range("N155").Value = 0
$ws.Range("O155").Value = 0
$ws.Range("P155").Value = 1
$ws.Range("Q155").Value = 0
$ws.Range("R155").Value = 0
$ws.Range("S155").Value = 0
$ws.Range("T155").Value = 0
$ws.Range("U155").Value = 0
$ws.Range("V155").Value = 0

# Row 156: Viktor Mikic
$ws.Range("A156").NumberFormat = "@"
$ws.Range("A156").Value = "2026-02-01"
$ws.Range("B156").Value = 'Undrafted'
$ws.Range("C156").Value = 'No'
$ws.Range("D156").Value = 'Viktor Mikic'
$ws.Range("E156").Value = 'FLA'
$ws.Range("F156").Value = 'ALA@FLA'
$ws.Range("G156").Value = '3:36 - 2nd Half'
$ws.Range("H156").Value = 0
$ws.Range("I156").Value = 0
$ws.Range("J156").Value = 0
$ws.Range("K156").Value = 0
$ws.Range("L156").Value = 0
$ws.Range("M156").Value = 0
$ws.Range("N156").Value = 0
$ws.Range("O156").Value = 0
$ws.Range("P156").Value = 1
$ws.Range("Q156").Value = 0
$ws.Range("R156").Value = 0
$ws.Range("S156").Value = 0
$ws.Range("T156").Value = 0
$ws.Range("U156").Value = 0
$ws.Range("V156").Value = 0

# ---- Part 4: OwnerTotals sheet totals ----
$ws2 = $wb.Worksheets.Item("OwnerTotals")
$ws2.Range("B2").Value = 99
$ws2.Range("B3").Value = 88
$ws2.Range("B4").Value = 82
$ws2.Range("B5").Value = 73
